# Update the LinkedIn carousel draft: swap the PFC Consulting / Andhra Pradesh
# transmission-project story for the Juniper Green Energy / Maharashtra solar
# story across all six slides (title + two body bullets on each slide).

$p = $ppt.ActivePresentation

$newTitle = "Juniper Green Energy Commissions Additional 72 MWp Solar Component of Hybrid Project in Maharashtra - Energetica India Magazine"

$bullets = @{
    1 = @(
        "Juniper Green Energy has commissioned an additional 72 MWp solar component.",
        "The project is part of a hybrid initiative located in Maharashtra."
    )
    2 = @(
        "The total capacity of the solar component is 72 MWp.",
        "This addition enhances the renewable energy output in Maharashtra."
    )
    3 = @(
        "The hybrid project aims to integrate solar energy with other renewable sources.",
        "Maharashtra is a key region for renewable energy development in India."
    )
    4 = @(
        "The commissioning of the solar component contributes to India's renewable energy targets.",
        "Juniper Green Energy is focused on expanding its renewable energy portfolio."
    )
    5 = @(
        "The project reflects ongoing investments in sustainable energy solutions.",
        "This initiative supports local energy needs and reduces carbon footprint."
    )
    6 = @(
        "The hybrid project is part of a broader strategy to enhance energy security.",
        "Juniper Green Energy's efforts align with national policies on renewable energy."
    )
}

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $s = $p.Slides.Item($i)

    # Shape 1: Title
    $titleShape = $s.Shapes.Item(1)
    $titleShape.TextFrame.TextRange.Runs(1).Text = $newTitle

    # Shape 2: Content placeholder - paragraph 1 is the image-missing
    # warning, paragraphs 2 and 3 are the two body bullets.
    $bodyShape = $s.Shapes.Item(2)
    $bodyText = $bodyShape.TextFrame.TextRange
    $pair = $bullets[$i]
    $bodyText.Paragraphs(2).Runs(1).Text = $pair[0]
    $bodyText.Paragraphs(3).Runs(1).Text = $pair[1]
}
